$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp banner.
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 13:16"

# Refresh the province case counts that came in with new data (this is what
# triggers the table, which is kept sorted descending by "Casos totales", to
# re-order itself below).
$ws.Range("A7").Value = "Araba/Alava"
$ws.Range("B7").Value = 1009
$ws.Range("C7").Value = 283
$ws.Range("D7").Value = 939
$ws.Range("E7").Value = 70

$ws.Range("A9").Value = "Bizkaia/Vizcaya"
$ws.Range("B9").Value = 1032
$ws.Range("C9").Value = 283
$ws.Range("D9").Value = 997
$ws.Range("E9").Value = 35

$ws.Range("A23").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B23").Value = 380
$ws.Range("C23").Value = 283
$ws.Range("D23").Value = 365
$ws.Range("E23").Value = 15

# Re-sort the full data table (rows 4-63) descending by column B (Casos
# totales), same as the source sheet does after every data refresh.
$dataRange = $ws.Range("A4:E63")
$sortKey = $ws.Range("B4")
$dataRange.Sort($sortKey, 2, $null, $null, $null, $null, $null, 0)
